# excelExport.xlsx edit script
# Implements: new defined-name/autofilter/dimension range A:V -> A:Z,
# four new columns (W:Z) with header / sub-header / placeholder rows,
# font rename ("MS Gothic" -> "SimSun") on the existing header fonts,
# and the "Normal" cell style Chinese display-name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) New header cells (row 2) - same style as column V (s="10")
# ---------------------------------------------------------------------
$ws.Range("V2:V4").Copy()
$ws.Range("W2:W4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("X2:X4").PasteSpecial(-4122)
$ws.Range("Y2:Y4").PasteSpecial(-4122)
$ws.Range("Z2:Z4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 headers --------------------------------------------------------
# W2: "考场号" - "场" "号" bold/宋体 (mirrors font used elsewhere in the sheet)
$ws.Range("W2").Value = "考场号"
$wFont = $ws.Range("W2").Characters(2, 2).Font
$wFont.Name = "宋体"
$wFont.Bold = $true
$wFont.ColorIndex = -4105

# X2: "座位号"
$ws.Range("X2").Value = "座位号"

# Y2: "准考证号" - "证" "号" bold/宋体
$ws.Range("Y2").Value = "准考证号"
$yFont = $ws.Range("Y2").Characters(3, 2).Font
$yFont.Name = "宋体"
$yFont.Bold = $true
$yFont.ColorIndex = -4105

# Z2: "准考证打印状态"
$ws.Range("Z2").Value = "准考证打印状态"

# ---------------------------------------------------------------------
# 2) New sub-header cells (row 3) - same style as column V (s="4")
# ---------------------------------------------------------------------
# W3: "${applyUsers.room.code}" - trailing "}" in Arial
$ws.Range("W3").Value = '${applyUsers.room.code}'
$w3Font = $ws.Range("W3").Characters(23, 1).Font
$w3Font.Name = "Arial"
$w3Font.Size = 10
$w3Font.ColorIndex = -4105

# X3: "${applyUsers.seat.code}"
$ws.Range("X3").Value = '${applyUsers.seat.code}'

# Y3: "${applyUsers.admission.code}"
$ws.Range("Y3").Value = '${applyUsers.admission.code}'

# Z3: ternary expression with mixed-font "是"/"否"
$z3Text = '${applyUsers.admission != null ? applyUsers.admission.printFlg ? "' + [char]0x662F + '":"' + [char]0x5426 + '":""}'
$ws.Range("Z3").Value = $z3Text
$yesFont = $ws.Range("Z3").Characters(67, 1).Font
$yesFont.Name = "ＭＳ Ｐゴシック"
$yesFont.Size = 10
$yesFont.ColorIndex = -4105
$mid1Font = $ws.Range("Z3").Characters(68, 3).Font
$mid1Font.Name = "Arial"
$mid1Font.Size = 10
$mid1Font.ColorIndex = -4105
$noFont = $ws.Range("Z3").Characters(71, 1).Font
$noFont.Name = "ＭＳ Ｐゴシック"
$noFont.Size = 10
$noFont.ColorIndex = -4105
$tailFont = $ws.Range("Z3").Characters(72, 5).Font
$tailFont.Name = "Arial"
$tailFont.Size = 10
$tailFont.ColorIndex = -4105

# ---------------------------------------------------------------------
# 3) New placeholder cells (row 4) - same content/style as V4
# ---------------------------------------------------------------------
$v4 = $ws.Range("V4").Value2
$ws.Range("W4").Value = $v4
$ws.Range("X4").Value = $v4
$ws.Range("Y4").Value = $v4
$ws.Range("Z4").Value = $v4

# ---------------------------------------------------------------------
# 4) Column widths for the new columns (match column V's width)
# ---------------------------------------------------------------------
$ws.Range("W1:Z1").EntireColumn.ColumnWidth = $ws.Columns.Item(22).ColumnWidth

# ---------------------------------------------------------------------
# 5) AutoFilter range + _FilterDatabase defined name A2:V2 -> A2:Z2
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A2:Z2").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$2:`$Z`$2"
    }
}

# ---------------------------------------------------------------------
# 6) Font rename: "ＭＳ ゴシック" -> "宋体" for the existing header fonts
# ---------------------------------------------------------------------
$ws.Range("A1").Font.Name = "宋体"
$ws.Range("D2:R2").Font.Name = "宋体"
$ws.Range("A2:C2").Font.Name = "宋体"
$ws.Range("S2:V2").Font.Name = "宋体"
$ws.Range("A3").Font.Name = "宋体"

# ---------------------------------------------------------------------
# 7) "Normal" cell style Chinese display name (best effort)
# ---------------------------------------------------------------------
$wb.Styles.Item(1).Name = "常规"
$null
